# Header updates for summer uploads
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (row 1) to reflect the new column headers.
$ws.Range("D1").Value = "Legal Given Name"
$ws.Range("E1").Value = "Birthdate"
$ws.Range("F1").Value = "Ministry Course Code and Level"
$ws.Range("H1").Value = "Final Percent"
$ws.Range("J1").Value = "Credits"

# Update the active selection on the sheet.
[void]$ws.Range("F1:J1").Select()
